# Zeitplanung Luminescence LED's - apply authored edit
# - Add new sheet "Tabelle1" between "Zeitplanung" and "Ist Arbeitszeit - Übersicht"
# - Fill in actual hours-worked data on the "Zeitplanung" sheet (columns C, J, K, P)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Tabelle1" worksheet right after "Zeitplanung"
# ---------------------------------------------------------------------------
$zeitplanung = $wb.Worksheets.Item("Zeitplanung")
$tabelle1 = $wb.Worksheets.Add($null, $zeitplanung)
$tabelle1.Name = "Tabelle1"

$tabelle1.Range("C3").Value = "Optimal"
$tabelle1.Range("C4").Value = "Normal"
$tabelle1.Range("F4").Formula = "=+(D3+2+D4+D5)/3"
$tabelle1.Range("C5").Value = "Worst Case"

# ---------------------------------------------------------------------------
# 2. Update the actual working hours ("Ist") entered into the Zeitplanung
#    schedule for the various tasks.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Zeitplanung")

# -- Administration, Planung ------------------------------------------------
# Arbeitsjournal nachführen (row 11)
$ws.Range("C11").Value = 3
$ws.Range("J11").Value = 0.5
$ws.Range("K11").Value = ""
$ws.Range("P11").Value = 0.5

# Zeitplanung erstellen (row 12)
$ws.Range("K12").Value = 0

# -- Analyse & Design ---------------------------------------------------------
# Anforderungsanalyse (Was) (row 15)
$ws.Range("C15").Value = 5
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 0
$ws.Range("P15").Value = 3

# -- Implementation -----------------------------------------------------------
# Anforderung #01 (row 19)
$ws.Range("C19").Value = 3
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 2

# Anforderung #02 (row 20)
$ws.Range("C20").Value = 10
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 6

# Anforderung #03 (row 21)
$ws.Range("C21").Value = 8

# Anforderung #04 (row 22)
$ws.Range("C22").Value = 6

# Anforderung #05 (row 23)
$ws.Range("C23").Value = 8

# Anforderung #06 (row 24)
$ws.Range("C24").Value = 4

# Anforderung #07 (row 25)
$ws.Range("C25").Value = 5

# Anforderung #08 (row 26)
$ws.Range("C26").Value = 10

# Anforderung #09 (row 27)
$ws.Range("C27").Value = 7

# Testfälle erstellen (row 28)
$ws.Range("C28").Value = 10

# Bugs fixen (row 29)
$ws.Range("C29").Value = 12

# -- Diverses -----------------------------------------------------------------
# Project Abstract (row 40)
$ws.Range("C40").Value = 4

# Meilenstein und ggf. Besprechung/Sitzung (row 41)
$ws.Range("C41").Value = 1.5
